# Weekly update: a new price record for "Feria Lagunitas de Puerto Montt - Ajo"
# is inserted as row 114 (pushing the existing rows 114-181 down to 115-182).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 114; this shifts rows 114..181 down to 115..182
# and grows the sheet dimension from R181 to R182 automatically.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new weekly record.
$ws.Range("A114").Value = 4
$ws.Range("B114").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C114").Value = "Los Lagos"
$ws.Range("D114").Value = 44529
$ws.Range("E114").Value = 10
$ws.Range("F114").Value = 100112003
$ws.Range("G114").Value = "Ajo"
$ws.Range("H114").Value = "Chino"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 60
$ws.Range("K114").Value = 21000
$ws.Range("L114").Value = 21000
$ws.Range("M114").Value = 21000
$ws.Range("N114").Value = '$/caja 10 kilos'
$ws.Range("O114").Value = "China"
$ws.Range("P114").Value = 2100
$ws.Range("Q114").Value = 10
$ws.Range("R114").Value = "Hortaliza"

# Note: Rows.Item(114).Insert() already carries the date-format style
# (the same one used by every other row's column D) onto the new D114,
# so no extra style assignment is needed here.
